$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41, pushing existing rows 41:94 down to 42:95
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with its data (mirrors the constant columns used
# throughout this block, plus the new record's specific values)
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value = "Ñuble"
$ws.Cells.Item(41, 4).Value = 44638
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 100112030
$ws.Cells.Item(41, 7).Value = "Poroto granado"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 60
$ws.Cells.Item(41, 11).Value = 20000
$ws.Cells.Item(41, 12).Value = 20000
$ws.Cells.Item(41, 13).Value = 20000
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(41, 16).Value = 800
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
